$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 113
$ws.Range("H113").Value = 3939.2
$ws.Range("I113").Value = 3166.6667
$ws.Range("J113").Value = 5098
$ws.Range("K113").Value = 3166.6667
$ws.Range("L113").Value = 5098
$ws.Range("M113").Value = 87.33329999999978
$ws.Range("N113").Value = -11606

# Row 121
$ws.Range("H121").Value = 1919.7
$ws.Range("J121").Value = 1919.7
$ws.Range("L121").Value = 5759.1
$ws.Range("N121").Value = -9253.1

# Row 132
$ws.Range("H132").Value = 1159.579
$ws.Range("I132").Value = 1207.8889
$ws.Range("J132").Value = 290
$ws.Range("K132").Value = 3623.6667
$ws.Range("L132").Value = 870
$ws.Range("M132").Value = -1093.6667
$ws.Range("N132").Value = -5930

# Row 137
$ws.Range("H137").Value = 4858.75
$ws.Range("J137").Value = 5840.273
$ws.Range("L137").Value = 17520.819
$ws.Range("N137").Value = -22620.819

# Row 138
$ws.Range("H138").Value = 2694.4314
$ws.Range("I138").Value = 1088.8
$ws.Range("J138").Value = 4988.1904
$ws.Range("K138").Value = 3266.4
$ws.Range("L138").Value = 14964.5712
$ws.Range("M138").Value = 1873.6
$ws.Range("N138").Value = -25244.5712

# Row 140
$ws.Range("H140").Value = 119449
$ws.Range("J140").Value = 119449
$ws.Range("L140").Value = 119449
$ws.Range("N140").Value = -129809

$ws = $wb.Worksheets.Item("ARM")
# Row 97
$ws.Range("H97").Value = 1149.6
$ws.Range("I97").Value = 506.5
$ws.Range("J97").Value = 3722
$ws.Range("K97").Value = 506.5
$ws.Range("L97").Value = 3722
$ws.Range("M97").Value = -10.5
$ws.Range("N97").Value = -4714

$ws = $wb.Worksheets.Item("BSM")
# Row 128
$ws.Range("H128").Value = 9580
$ws.Range("I128").Value = 9580
$ws.Range("K128").Value = 28740
$ws.Range("M128").Value = -26250

# Row 132
$ws.Range("H132").Value = 114674.25
$ws.Range("I132").Value = 0
$ws.Range("J132").Value = 114674.25
$ws.Range("K132").Value = 0
$ws.Range("L132").Value = 114674.25
$ws.Range("M132").ClearContents()
$ws.Range("N132").Value = -124794.25

$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Range("H31").Value = 2285.375
$ws.Range("I31").Value = 2483.7144
$ws.Range("J31").Value = 897
$ws.Range("K31").Value = 2483.7144
$ws.Range("L31").Value = 897
$ws.Range("M31").Value = -2188.7144
$ws.Range("N31").Value = -1487

# Row 34
$ws.Range("H34").Value = 2285.375
$ws.Range("I34").Value = 2483.7144
$ws.Range("J34").Value = 897
$ws.Range("K34").Value = 2483.7144
$ws.Range("L34").Value = 897
$ws.Range("M34").Value = -2281.7144
$ws.Range("N34").Value = -1301

# Row 86
$ws.Range("H86").Value = 19182.777
$ws.Range("I86").Value = 4906
$ws.Range("J86").Value = 41617.715
$ws.Range("K86").Value = 4906
$ws.Range("L86").Value = 41617.715
$ws.Range("M86").Value = -3783
$ws.Range("N86").Value = -43863.715

# Row 89
$ws.Range("H89").Value = 19182.777
$ws.Range("I89").Value = 4906
$ws.Range("J89").Value = 41617.715
$ws.Range("K89").Value = 24530
$ws.Range("L89").Value = 208088.575
$ws.Range("M89").Value = -18914
$ws.Range("N89").Value = -219320.575

# Row 99
$ws.Range("H99").Value = 4247.25
$ws.Range("I99").Value = 4333
$ws.Range("J99").Value = 3990
$ws.Range("K99").Value = 4333
$ws.Range("L99").Value = 3990
$ws.Range("M99").Value = -2835
$ws.Range("N99").Value = -6986

# Row 126
$ws.Range("H126").Value = 4247.25
$ws.Range("I126").Value = 4333
$ws.Range("J126").Value = 3990
$ws.Range("K126").Value = 12999
$ws.Range("L126").Value = 11970
$ws.Range("M126").Value = -10529
$ws.Range("N126").Value = -16910

# Row 134
$ws.Range("H134").Value = 5732.125
$ws.Range("I134").Value = 5732.125
$ws.Range("K134").Value = 17196.375
$ws.Range("M134").Value = -14661.375

$ws = $wb.Worksheets.Item("CUL")
# Row 137
$ws.Range("H137").Value = 0
$ws.Range("I137").Value = 0
$ws.Range("J137").Value = 0
$ws.Range("K137").Value = 0
$ws.Range("L137").Value = 0
$ws.Range("M137").ClearContents()
$ws.Range("N137").ClearContents()

$ws = $wb.Worksheets.Item("GSM")
# Row 80
$ws.Range("H80").Value = 4022.25
$ws.Range("I80").Value = 3478.2727
$ws.Range("K80").Value = 3478.2727
$ws.Range("M80").Value = -2480.2727

# Row 83
$ws.Range("H83").Value = 4022.25
$ws.Range("I83").Value = 3478.2727
$ws.Range("K83").Value = 17391.3635
$ws.Range("M83").Value = -12399.3635

# Row 126
$ws.Range("H126").Value = 0
$ws.Range("I126").Value = 0
$ws.Range("K126").Value = 0
$ws.Range("M126").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
# Row 17
$ws.Range("H17").Value = 20000
$ws.Range("J17").Value = 20000
$ws.Range("L17").Value = 20000
$ws.Range("N17").Value = -20340

# Row 82
$ws.Range("H82").Value = 1799.5
$ws.Range("I82").Value = 1799.5
$ws.Range("J82").Value = 1799.5
$ws.Range("K82").Value = 1799.5
$ws.Range("L82").Value = 1799.5
$ws.Range("M82").Value = -1438.5
$ws.Range("N82").Value = -2521.5

# Row 85
$ws.Range("H85").Value = 1799.5
$ws.Range("I85").Value = 1799.5
$ws.Range("J85").Value = 1799.5
$ws.Range("K85").Value = 1799.5
$ws.Range("L85").Value = 1799.5
$ws.Range("M85").Value = -551.5
$ws.Range("N85").Value = -4295.5

$ws = $wb.Worksheets.Item("WVR")
# Row 14
$ws.Range("J14").Value = 0
$ws.Range("L14").Value = 0
$ws.Range("N14").ClearContents()

# Row 62
$ws.Range("H62").Value = 30000
$ws.Range("J62").Value = 30000
$ws.Range("L62").Value = 30000
$ws.Range("N62").Value = -31248

# Row 65
$ws.Range("H65").Value = 30000
$ws.Range("J65").Value = 30000
$ws.Range("L65").Value = 150000
$ws.Range("N65").Value = -156240

# Row 81
$ws.Range("H81").Value = 8704.75
$ws.Range("I81").Value = 1326.8
$ws.Range("K81").Value = 2653.6
$ws.Range("M81").Value = -1592.6

# Row 84
$ws.Range("H84").Value = 8704.75
$ws.Range("I84").Value = 1326.8
$ws.Range("K84").Value = 13268
$ws.Range("M84").Value = -7964

# Row 132
$ws.Range("H132").Value = 1688.4
$ws.Range("I132").Value = 1593.3334
$ws.Range("K132").Value = 4780.0002
$ws.Range("M132").Value = -2250.0002

# Row 136
$ws.Range("H136").Value = 1126.1154
$ws.Range("I136").Value = 947.86365
$ws.Range("K136").Value = 2843.59095
$ws.Range("M136").Value = -293.5909499999998
